# Refresh the cryptos worksheet: update Price (col D) / Volume(1h) (col E)
# figures for each coin row to the latest scrape, and re-rank rows 17/18
# (WrappedEther now outranks TRON), per the scheduled GitHub Actions job.
#
# Values are written indirectly through a scratch cell that is forced to
# text format, then copied with "paste values only" onto the real target.
# This keeps every touched cell's stored type as plain text (matching the
# sheet's original inlineStr cells) even for values that look like plain
# numbers (e.g. "558.76"), instead of letting Excel's normal Range.Value
# auto-detection silently convert them into numeric cells. The scratch
# cell is fully cleared (content + formatting) afterwards so the sheet's
# used range / styles are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
2|D|65.436.78
2|E|  -0.31%  
3|D|3.376.51
3|E|  -1.37%  
4|E|  +0.07%  
5|D|558.76
5|E|  -0.85%  
6|D|174.81
6|E|  -0.93%  
7|D|0.628
7|E|  +0.28%  
8|D|3.365.29
8|E|  -1.47%  
9|E|  +0.07%  
10|D|0.173
10|E|  +1.28%  
11|D|0.632
11|E|  -0.33%  
12|D|53.33
12|E|  -3.33%  
13|D|0.0000276
13|E|  -1.63%  
14|D|9.16
14|E|  -0.16%  
15|D|3.932.99
15|E|  -0.85%  
16|D|18.19
16|E|  -1.26%  
17|B|WrappedEther
17|C|https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth
17|D|3.389.14
17|E|  -0.74%  
18|B|TRON
18|C|https://coinranking.com/coin/qUhEFk1I61atv+tron-trx
18|D|0.119
18|E|  +0.16%  
19|D|65.594.06
19|E|  -0.04%  
20|D|11.80
20|E|  -1.51%  
21|D|0.995
21|E|  -0.41%  
22|D|480.82
22|E|  +2.66%  
23|D|4.98
23|E|  -0.43%  
24|D|90.03
24|E|  +3.80%  
25|D|4.09
25|E|  -1.70%  
26|D|14.23
26|E|  +3.65%  
27|D|2.90
27|E|  -0.21%  
28|D|10.58
28|E|  -3.28%  
29|D|8.69
29|E|  -2.47%  
30|D|31.18
30|E|  +1.04%  
31|D|6.54
31|E|  -2.73%  
32|D|63.64
32|E|  +5.13%  
33|D|11.40
33|E|  -1.61%  
34|D|572.32
34|E|  -2.67%  
35|E|  -1.53%  
36|E|  +0.00%  
37|D|3.61
37|E|  +2.96%  
38|E|  +0.10%  
39|D|35.64
39|E|  -1.37%  
40|D|0.373
40|E|  -0.81%  
41|D|0.0₃0737
41|E|  -3.23%  
42|D|3.100.39
42|E|  -0.42%  
43|E|  -3.07%  
44|D|0.0414
44|E|  -0.58%  
45|D|0.133
45|E|  -0.65%  
46|E|  -2.03%  
47|D|2.43
47|E|  -3.65%  
48|E|  +0.11%  
49|D|140.61
49|E|  +2.81%  
50|D|2.58
50|E|  -0.01%  
51|D|8.39
51|E|  -0.05%  
"@

$lines = $data -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }

$helper = $ws.Range("Z1")

foreach ($line in $lines) {
    $parts = $line.Split('|', 3)
    $row = [int]$parts[0]
    $col = $parts[1]
    $value = $parts[2]

    $helper.NumberFormat = "@"
    $helper.Value = $value
    $helper.Copy()

    $target = $ws.Range("$col$row")
    $target.PasteSpecial(-4163)  # xlPasteValues

    $helper.Clear()
}
